$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C17").Value = "Use of contraceptive pills or injections"
$ws.Range("C24").Value = "History of diabetes"
$ws.Range("C35").Value = "Screening, skin cancer"
$ws.Range("C36").Value = "Screening, mammography"
$ws.Range("C37").Value = "Screening cervical, smear test"
$ws.Range("C59").Value = "Type of Cancer (ICD 10, 3 digits,e.g. C18)"
$ws.Range("C66").Value = "Body Mass Index at baseline"
$ws.Range("C67").Value = "Body Mass Index at follow-up"
$ws.Range("C68").Value = "Body Mass Index Standard Deviation Score at baseline (children studies)"
$ws.Range("C69").Value = "Body Mass Index Standard Deviation Score at follow-up (children studies)"
$ws.Range("C76").Value = "Body fat precent at follow-up"
$ws.Range("C77").Value = "Body fat precent at baseline"
$ws.Range("C94").Value = "Daily glycaemic load"
$ws.Range("C98").Value = "Intake of cakes and fine bakery products [g/d]"
$ws.Range("C99").Value = "Intake of fruit and vegetable juices [g/d]"
$ws.Range("C100").Value = "Intake of soft drinks [g/d]"
$ws.Range("C103").Value = "Total legumes intake [g/d]"
$ws.Range("C104").Value = "Total fruit intake [g/d]"

$ws.Range("C2:C109").Select() | Out-Null
